# Weekly driver report update for 2025-04-21
# Updates Critical Minutes (C) and Good Roaming Calculation (%) (D) values
# for the "Bad Drivers" table, plus the Totals row (C6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wi-Fi 6E AX210 160MHz - 23.60.1.2
$ws.Range("C3").Value = 627
$ws.Range("D3").Value = 87.40000000000001

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4
$ws.Range("C5").Value = 266
$ws.Range("D5").Value = 98.8

# Row 6: Totals
$ws.Range("C6").Value = 1198
